$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.692.52"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.633.47"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.22"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.22"
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0613"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("E11").Value = "  -3.06%  "
$ws.Range("D12").Value = "1.863.72"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").Value = "1.634.93"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.553"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.20"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "27.656.60"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.22"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "0.0₃0721"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("E20").Value = "  -1.68%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.66"
$ws.Range("E22").Value = "  +4.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.35"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.14"
$ws.Range("E24").Value = "  +2.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.82"
$ws.Range("E25").Value = "  -1.46%  "
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.61"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0482"
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("D33").Value = "1.472.81"
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  -1.62%  "
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("E37").Value = "  +6.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.878"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "67.91"
$ws.Range("E42").Value = "  -1.63%  "
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.35"
$ws.Range("E45").Value = "  -4.68%  "
$ws.Range("D46").Value = "1.773.56"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.53"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.73"
$ws.Range("E51").Value = "  -1.32%  "
